$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the header row two columns to the left (C3:I3 -> A3:G3),
# bringing values + full formatting (fills, borders, fonts, alignment) with it.
$ws.Range("C3:I3").Copy($ws.Range("A3"))

# Drop the now-stale copies left behind in the old H3:I3 positions.
$ws.Range("H3:I3").Clear()

# New 8th column: "rIC" header with its own accent fill.
$ws.Range("H3").Value = "rIC"
$ws.Range("H3").Font.Bold = $true
$ws.Range("H3").Font.Size = 12
$ws.Range("H3").HorizontalAlignment = -4108
$ws.Range("H3").VerticalAlignment = -4108
$ws.Range("H3").Interior.Pattern = 1
$ws.Range("H3").Interior.ThemeColor = 2
$ws.Range("H3").Interior.TintAndShade = -0.249977111117893

# Column widths now span A:H instead of C:I.
$ws.Range("A:H").ColumnWidth = 44.7109375

# Match the saved selection from the edit.
$ws.Range("D7").Select()
